$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newDates = @(
    "Fri Oct 06 11:15:46 EDT 2023",
    "Fri Oct 06 11:15:56 EDT 2023",
    "Fri Oct 06 11:16:06 EDT 2023",
    "Fri Oct 06 11:16:17 EDT 2023",
    "Fri Oct 06 11:16:27 EDT 2023",
    "Fri Oct 06 11:16:37 EDT 2023",
    "Fri Oct 06 11:16:48 EDT 2023",
    "Fri Oct 06 11:16:58 EDT 2023",
    "Fri Oct 06 11:17:08 EDT 2023",
    "Fri Oct 06 11:17:19 EDT 2023",
    "Fri Oct 06 11:17:29 EDT 2023",
    "Fri Oct 06 11:17:40 EDT 2023",
    "Fri Oct 06 11:17:50 EDT 2023",
    "Fri Oct 06 11:18:00 EDT 2023",
    "Fri Oct 06 11:18:10 EDT 2023",
    "Fri Oct 06 11:18:20 EDT 2023",
    "Fri Oct 06 11:18:30 EDT 2023",
    "Fri Oct 06 11:18:41 EDT 2023",
    "Fri Oct 06 11:18:51 EDT 2023",
    "Fri Oct 06 11:19:01 EDT 2023",
    "Fri Oct 06 11:19:11 EDT 2023",
    "Fri Oct 06 11:19:22 EDT 2023",
    "Fri Oct 06 11:19:32 EDT 2023",
    "Fri Oct 06 11:19:42 EDT 2023",
    "Fri Oct 06 11:19:53 EDT 2023",
    "Fri Oct 06 11:20:03 EDT 2023",
    "Fri Oct 06 11:20:14 EDT 2023",
    "Fri Oct 06 11:20:24 EDT 2023"
)

$startRow = 2
for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("B$row").Value = $newDates[$i]
}
